$p = $ppt.ActivePresentation

# Add a new slide using the "Title Slide" layout (layout index 1 == ppLayoutTitle)
$layout = $p.SlideMaster.CustomLayouts.Item(1)
$s = $p.Slides.AddSlide(1, $layout)

# Set the title text
$s.Shapes.Item(1).TextFrame.TextRange.Text = "TEST"
